$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 13 to make room for the split "Docentes responsaveis" rows
# (this shifts existing rows 13-23 down to 15-25, and their content/formatting travels with them)
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()

# Update cell contents that are new or changed
$ws.Range('B10').Value = 'Conduzir os alunos no desenvolvimento de um projeto de conclusão de curso, sobre um tema específico relacionado à Engenharia Bioquímica.'
$ws.Range('C10').Value = 'Conduzir os alunos no desenvolvimento de um projeto de conclusão de curso, sobre um tema específico relacionado à Engenharia Bioquímica.'
$ws.Range('B13').Value = '1304060 - Maria das Graças de Almeida Felipe'
$ws.Range('C13').Value = '1304060 - Maria das Graças de Almeida Felipe'
$ws.Range('B14').Value = '8853480 - Tatiane da Franca Silva'
$ws.Range('C14').Value = '8853480 - Tatiane da Franca Silva'
$ws.Range('B15').Value = '1) Metodologia de pesquisa em Engenharia 2) Elementos de projeto de monografia 3) Métodos de pesquisa4) Normas de um projeto de pesquisa.5) Etapas de elaboração de textos científicos'
$ws.Range('C15').Value = '1) Metodologia de pesquisa em Engenharia 2) Elementos de projeto de monografia 3) Métodos de pesquisa4) Normas de um projeto de pesquisa.5) Etapas de elaboração de textos científicos'
$ws.Range('B17').Value = '1. Metodologia Cientifica em engenharia: Concepção e definição. 2. Monografia Cientifica: O que é um projeto de pesquisa. Os elementos que compõem um projeto de pesquisa. 3. Métodos de Pesquisa utilizados na Engenharia Bioquímica. 4. Normas para elaboração de do texto e das Referências Bibliográficas.5. Etapas de elaboração de textos científicos'
$ws.Range('C17').Value = '1. Metodologia Cientifica em engenharia: Concepção e definição. 2. Monografia Cientifica: O que é um projeto de pesquisa. Os elementos que compõem um projeto de pesquisa. 3. Métodos de Pesquisa utilizados na Engenharia Bioquímica. 4. Normas para elaboração de do texto e das Referências Bibliográficas.5. Etapas de elaboração de textos científicos'
$ws.Range('B20').Value = 'Apresentação de um pré-projeto e um projeto. O projeto será avaliado por dois examinadores. A média (M) será calculada levando-se a nota do pré-projeto(NPP) e a média da nova dos dois examinadores (NP) M = (0,3PP + 0,7NP), conforme Norma para Trabalho de Conclusão de Curso do curso de Engenharia Bioquímica.'
$ws.Range('C20').Value = 'Apresentação de um pré-projeto e um projeto. O projeto será avaliado por dois examinadores. A média (M) será calculada levando-se a nota do pré-projeto(NPP) e a média da nova dos dois examinadores (NP) M = (0,3PP + 0,7NP), conforme Norma para Trabalho de Conclusão de Curso do curso de Engenharia Bioquímica.'
$ws.Range('B21').Value = 'M=≥ 5,0 para ser aprovado'
$ws.Range('C21').Value = 'M=≥ 5,0 para ser aprovado'
$ws.Range('B22').Value = '(NF+RP)/2 ≥ 5,0 para ser aprovado, onde RP é a nota do projeto modificado apresentado.'
$ws.Range('C22').Value = '(NF+RP)/2 ≥ 5,0 para ser aprovado, onde RP é a nota do projeto modificado apresentado.'
$ws.Range('B23').Value = 'SANTOS, C. R. Trabalho de Conclusão de Curso – Guia de elaboração passo a passo, Cengage Learning, 2010.ANDRADE, Maria Margarida de. Introdução à metodologia do trabalho científico. 10.ed. São Paulo: Atlas, 2010.BOOTH, W.; COLOMB, G.; WILLIAMS, J. A arte da Pesquisa. 3 ed. Martins Fontes. São Paulo. 2005. GIL, A.C. Como elaborar projetos de pesquisa. 5 ed. Atlas, São Paulo, 2010. MEDEIROS, J. B. Redação Cientifica: A Prática de Fichamentos, Resumos e Resenhas. 11 ed. São Paulo: Atlas, 2009SERAFINI, Maria José. Como escrever textos. 5.ed. São Paulo: Globo, 1992.SEVERINO, Antonio Joaquim. Metodologia do trabalho científico. 23.ed. São Paulo: Cortez, 2009.'
$ws.Range('C23').Value = 'SANTOS, C. R. Trabalho de Conclusão de Curso – Guia de elaboração passo a passo, Cengage Learning, 2010.ANDRADE, Maria Margarida de. Introdução à metodologia do trabalho científico. 10.ed. São Paulo: Atlas, 2010.BOOTH, W.; COLOMB, G.; WILLIAMS, J. A arte da Pesquisa. 3 ed. Martins Fontes. São Paulo. 2005. GIL, A.C. Como elaborar projetos de pesquisa. 5 ed. Atlas, São Paulo, 2010. MEDEIROS, J. B. Redação Cientifica: A Prática de Fichamentos, Resumos e Resenhas. 11 ed. São Paulo: Atlas, 2009SERAFINI, Maria José. Como escrever textos. 5.ed. São Paulo: Globo, 1992.SEVERINO, Antonio Joaquim. Metodologia do trabalho científico. 23.ed. São Paulo: Cortez, 2009.'

# Narrow column A definition so it no longer overlaps column B (cosmetic range cleanup)
$ws.Columns.Item(1).ColumnWidth = 30.7109375
